$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value2 = 1468.7778
$ws.Cells.Item(2,10).Value2 = 2633.3333
$ws.Cells.Item(2,12).Value2 = 2633.3333
$ws.Cells.Item(2,14).Value2 = -2859.3333

$ws.Cells.Item(40,8).Value2 = 3430.2
$ws.Cells.Item(40,9).Value2 = 2750
$ws.Cells.Item(40,10).Value2 = 3534.8462
$ws.Cells.Item(40,11).Value2 = 2750
$ws.Cells.Item(40,12).Value2 = 3534.8462
$ws.Cells.Item(40,13).Value2 = -2575
$ws.Cells.Item(40,14).Value2 = -3884.8462

$ws.Cells.Item(41,8).Value2 = 1234
$ws.Cells.Item(41,9).Value2 = 900
$ws.Cells.Item(41,11).Value2 = 900
$ws.Cells.Item(41,13).Value2 = -460

$ws.Cells.Item(70,8).Value2 = 1269094.5
$ws.Cells.Item(70,9).Value2 = 3376171.2
$ws.Cells.Item(70,10).Value2 = 4848.6
$ws.Cells.Item(70,11).Value2 = 10128513.6
$ws.Cells.Item(70,12).Value2 = 14545.8
$ws.Cells.Item(70,13).Value2 = -10128243.6
$ws.Cells.Item(70,14).Value2 = -15085.8

$ws.Cells.Item(73,8).Value2 = 1269094.5
$ws.Cells.Item(73,9).Value2 = 3376171.2
$ws.Cells.Item(73,10).Value2 = 4848.6
$ws.Cells.Item(73,11).Value2 = 10128513.6
$ws.Cells.Item(73,12).Value2 = 14545.8
$ws.Cells.Item(73,13).Value2 = -10127577.6
$ws.Cells.Item(73,14).Value2 = -16417.8

$ws.Cells.Item(80,8).Value2 = 1666.6666
$ws.Cells.Item(80,9).Value2 = 1000
$ws.Cells.Item(80,10).Value2 = 2000
$ws.Cells.Item(80,11).Value2 = 3000
$ws.Cells.Item(80,12).Value2 = 6000
$ws.Cells.Item(80,13).Value2 = -2002
$ws.Cells.Item(80,14).Value2 = -7996

$ws.Cells.Item(83,8).Value2 = 1666.6666
$ws.Cells.Item(83,9).Value2 = 1000
$ws.Cells.Item(83,10).Value2 = 2000
$ws.Cells.Item(83,11).Value2 = 9000
$ws.Cells.Item(83,12).Value2 = 18000
$ws.Cells.Item(83,13).Value2 = -4008
$ws.Cells.Item(83,14).Value2 = -27984

$ws.Cells.Item(86,8).Value2 = 16742.334
$ws.Cells.Item(86,9).Value2 = 3000
$ws.Cells.Item(86,10).Value2 = 19490.8
$ws.Cells.Item(86,11).Value2 = 3000
$ws.Cells.Item(86,12).Value2 = 19490.8
$ws.Cells.Item(86,13).Value2 = -1877
$ws.Cells.Item(86,14).Value2 = -21736.8

$ws.Cells.Item(89,8).Value2 = 16742.334
$ws.Cells.Item(89,9).Value2 = 3000
$ws.Cells.Item(89,10).Value2 = 19490.8
$ws.Cells.Item(89,11).Value2 = 15000
$ws.Cells.Item(89,12).Value2 = 97454
$ws.Cells.Item(89,13).Value2 = -9384
$ws.Cells.Item(89,14).Value2 = -108686

$ws.Cells.Item(112,8).Value2 = 2645.5715
$ws.Cells.Item(112,9).Value2 = 1162
$ws.Cells.Item(112,10).Value2 = 3239
$ws.Cells.Item(112,11).Value2 = 3486
$ws.Cells.Item(112,12).Value2 = 9717
$ws.Cells.Item(112,13).Value2 = -2378
$ws.Cells.Item(112,14).Value2 = -11933

$ws.Cells.Item(113,8).Value2 = 6859.6
$ws.Cells.Item(113,9).Value2 = 4112
$ws.Cells.Item(113,10).Value2 = 9999.714
$ws.Cells.Item(113,11).Value2 = 4112
$ws.Cells.Item(113,12).Value2 = 9999.714
$ws.Cells.Item(113,13).Value2 = -858
$ws.Cells.Item(113,14).Value2 = -16507.714

$ws.Cells.Item(115,8).Value2 = 243
$ws.Cells.Item(115,9).Value2 = 243
$ws.Cells.Item(115,11).Value2 = 729
$ws.Cells.Item(115,13).Value2 = 838

$ws.Cells.Item(116,8).Value2 = 6234.3887
$ws.Cells.Item(116,9).Value2 = 4970
$ws.Cells.Item(116,10).Value2 = 6866.5835
$ws.Cells.Item(116,11).Value2 = 4970
$ws.Cells.Item(116,12).Value2 = 6866.5835
$ws.Cells.Item(116,13).Value2 = -1528
$ws.Cells.Item(116,14).Value2 = -13750.5835

$ws.Cells.Item(129,8).Value2 = 3546.9375
$ws.Cells.Item(129,9).Value2 = 1262.6666
$ws.Cells.Item(129,10).Value2 = 4074.077
$ws.Cells.Item(129,11).Value2 = 3787.9998
$ws.Cells.Item(129,12).Value2 = 12222.231
$ws.Cells.Item(129,13).Value2 = 1212.0002
$ws.Cells.Item(129,14).Value2 = -22222.231

$ws.Cells.Item(132,8).Value2 = 33337344
$ws.Cells.Item(132,9).Value2 = 41670804
$ws.Cells.Item(132,11).Value2 = 125012412
$ws.Cells.Item(132,13).Value2 = -125009882

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(21,8).Value2 = 13082
$ws.Cells.Item(21,9).Value2 = 8000
$ws.Cells.Item(21,10).Value2 = 14352.5
$ws.Cells.Item(21,11).Value2 = 8000
$ws.Cells.Item(21,12).Value2 = 14352.5
$ws.Cells.Item(21,13).Value2 = -7626
$ws.Cells.Item(21,14).Value2 = -15100.5

$ws.Cells.Item(32,8).Value2 = 7037.636
$ws.Cells.Item(32,9).Value2 = 7037.636
$ws.Cells.Item(32,11).Value2 = 7037.636
$ws.Cells.Item(32,13).Value2 = -6750.636

$ws.Cells.Item(80,8).Value2 = 37500

$ws.Cells.Item(83,8).Value2 = 37500

$ws.Cells.Item(102,8).Value2 = 744
$ws.Cells.Item(102,9).Value2 = 741.6
$ws.Cells.Item(102,10).Value2 = 750
$ws.Cells.Item(102,11).Value2 = 741.6
$ws.Cells.Item(102,12).Value2 = 750
$ws.Cells.Item(102,13).Value2 = 880.4
$ws.Cells.Item(102,14).Value2 = -3994

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82,8).Value2 = 16768.111

$ws.Cells.Item(85,8).Value2 = 16768.111

$ws.Cells.Item(86,8).Value2 = 1683.3334
$ws.Cells.Item(86,9).Value2 = 1821
$ws.Cells.Item(86,10).Value2 = 995
$ws.Cells.Item(86,11).Value2 = 1821
$ws.Cells.Item(86,12).Value2 = 995
$ws.Cells.Item(86,13).Value2 = -698
$ws.Cells.Item(86,14).Value2 = -3241

$ws.Cells.Item(89,8).Value2 = 1683.3334
$ws.Cells.Item(89,9).Value2 = 1821
$ws.Cells.Item(89,10).Value2 = 995
$ws.Cells.Item(89,11).Value2 = 9105
$ws.Cells.Item(89,12).Value2 = 4975
$ws.Cells.Item(89,13).Value2 = -3489
$ws.Cells.Item(89,14).Value2 = -16207

$ws.Cells.Item(94,8).Value2 = 3749.75
$ws.Cells.Item(94,9).Value2 = 3749.75
$ws.Cells.Item(94,11).Value2 = 3749.75
$ws.Cells.Item(94,13).Value2 = -3298.75

$ws.Cells.Item(99,8).Value2 = 800
$ws.Cells.Item(99,10).Value2 = 0
$ws.Cells.Item(99,12).Value2 = 0
$ws.Cells.Item(99,14).ClearContents()

$ws.Cells.Item(105,8).Value2 = 3850
$ws.Cells.Item(105,9).Value2 = 0
$ws.Cells.Item(105,10).Value2 = 3850
$ws.Cells.Item(105,11).Value2 = 0
$ws.Cells.Item(105,12).Value2 = 3850
$ws.Cells.Item(105,13).ClearContents()
$ws.Cells.Item(105,14).Value2 = -7344

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value2 = 4245
$ws.Cells.Item(31,9).Value2 = 2733
$ws.Cells.Item(31,10).Value2 = 5152.2
$ws.Cells.Item(31,11).Value2 = 2733
$ws.Cells.Item(31,12).Value2 = 5152.2
$ws.Cells.Item(31,13).Value2 = -2438
$ws.Cells.Item(31,14).Value2 = -5742.2

$ws.Cells.Item(34,8).Value2 = 4245
$ws.Cells.Item(34,9).Value2 = 2733
$ws.Cells.Item(34,10).Value2 = 5152.2
$ws.Cells.Item(34,11).Value2 = 2733
$ws.Cells.Item(34,12).Value2 = 5152.2
$ws.Cells.Item(34,13).Value2 = -2531
$ws.Cells.Item(34,14).Value2 = -5556.2

$ws.Cells.Item(62,8).Value2 = 6766.3335
$ws.Cells.Item(62,9).Value2 = 6987.125
$ws.Cells.Item(62,10).Value2 = 5000
$ws.Cells.Item(62,11).Value2 = 6987.125
$ws.Cells.Item(62,12).Value2 = 5000
$ws.Cells.Item(62,13).Value2 = -6363.125
$ws.Cells.Item(62,14).Value2 = -6248

$ws.Cells.Item(65,8).Value2 = 6766.3335
$ws.Cells.Item(65,9).Value2 = 6987.125
$ws.Cells.Item(65,10).Value2 = 5000
$ws.Cells.Item(65,11).Value2 = 34935.625
$ws.Cells.Item(65,12).Value2 = 25000
$ws.Cells.Item(65,13).Value2 = -31815.625
$ws.Cells.Item(65,14).Value2 = -31240

$ws.Cells.Item(105,8).Value2 = 1322.625
$ws.Cells.Item(105,9).Value2 = 1317.4667
$ws.Cells.Item(105,11).Value2 = 1317.4667
$ws.Cells.Item(105,13).Value2 = 429.5333000000001

$ws.Cells.Item(107,8).Value2 = 448.25
$ws.Cells.Item(107,9).Value2 = 459.83334
$ws.Cells.Item(107,10).Value2 = 413.5
$ws.Cells.Item(107,11).Value2 = 459.83334
$ws.Cells.Item(107,12).Value2 = 413.5
$ws.Cells.Item(107,13).Value2 = 1460.16666
$ws.Cells.Item(107,14).Value2 = -4253.5

$ws.Cells.Item(132,8).Value2 = 1998.2
$ws.Cells.Item(132,9).Value2 = 1998.2
$ws.Cells.Item(132,11).Value2 = 5994.6
$ws.Cells.Item(132,13).Value2 = -3464.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(47,8).Value2 = 84.5
$ws.Cells.Item(47,9).Value2 = 84.5
$ws.Cells.Item(47,11).Value2 = 253.5
$ws.Cells.Item(47,13).Value2 = 177.5

$ws.Cells.Item(59,8).Value2 = 0
$ws.Cells.Item(59,10).Value2 = 0
$ws.Cells.Item(59,12).Value2 = 0
$ws.Cells.Item(59,14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3,8).Value2 = 253575.5
$ws.Cells.Item(3,10).Value2 = 0
$ws.Cells.Item(3,12).Value2 = 0
$ws.Cells.Item(3,14).ClearContents()

$ws.Cells.Item(46,8).Value2 = 13142.571

$ws.Cells.Item(80,8).Value2 = 7934.4287
$ws.Cells.Item(80,9).Value2 = 7108.2
$ws.Cells.Item(80,11).Value2 = 7108.2
$ws.Cells.Item(80,13).Value2 = -6110.2

$ws.Cells.Item(83,8).Value2 = 7934.4287
$ws.Cells.Item(83,9).Value2 = 7108.2
$ws.Cells.Item(83,11).Value2 = 35541
$ws.Cells.Item(83,13).Value2 = -30549

$ws.Cells.Item(122,8).Value2 = 44765
$ws.Cells.Item(122,9).Value2 = 49853.332
$ws.Cells.Item(122,11).Value2 = 149559.996
$ws.Cells.Item(122,13).Value2 = -147109.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16,8).Value2 = 2784.4614
$ws.Cells.Item(16,9).Value2 = 3249.8
$ws.Cells.Item(16,10).Value2 = 1233.3334
$ws.Cells.Item(16,11).Value2 = 3249.8
$ws.Cells.Item(16,12).Value2 = 1233.3334
$ws.Cells.Item(16,13).Value2 = -3079.8
$ws.Cells.Item(16,14).Value2 = -1573.3334

$ws.Cells.Item(46,8).Value2 = 2805.7058
$ws.Cells.Item(46,9).Value2 = 2634.5454
$ws.Cells.Item(46,10).Value2 = 3119.5
$ws.Cells.Item(46,11).Value2 = 2634.5454
$ws.Cells.Item(46,12).Value2 = 3119.5
$ws.Cells.Item(46,13).Value2 = -2446.5454
$ws.Cells.Item(46,14).Value2 = -3495.5

$ws.Cells.Item(136,8).Value2 = 3502.6667
$ws.Cells.Item(136,9).Value2 = 3502.6667
$ws.Cells.Item(136,11).Value2 = 10508.0001
$ws.Cells.Item(136,13).Value2 = -7958.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64,8).Value2 = 49999
$ws.Cells.Item(64,10).Value2 = 49999
$ws.Cells.Item(64,12).Value2 = 49999
$ws.Cells.Item(64,14).Value2 = -50495

$ws.Cells.Item(67,8).Value2 = 49999
$ws.Cells.Item(67,10).Value2 = 49999
$ws.Cells.Item(67,12).Value2 = 49999
$ws.Cells.Item(67,14).Value2 = -51715
